$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the affected columns so numeric-looking strings
# (e.g. "288.65", "0.83%", "21") are stored as text, matching the
# original inlineStr text cells rather than being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "288.65"
$ws.Range("E2").Value = "0.83%"
$ws.Range("G2").Value = "21"

$ws.Range("D3").Value = "29.25"
$ws.Range("E3").Value = "1.66%"
$ws.Range("G3").Value = "21"

$ws.Range("D4").Value = "5.293"
$ws.Range("E4").Value = "4.70%"
$ws.Range("G4").Value = "21"

$ws.Range("D5").Value = "0.07046"
$ws.Range("E5").Value = "4.98%"
$ws.Range("G5").Value = "21"

$ws.Range("D6").Value = "7.449"
$ws.Range("E6").Value = "1.66%"
$ws.Range("G6").Value = "21"

$ws.Range("D7").Value = "3.575"
$ws.Range("E7").Value = "5.42%"
$ws.Range("G7").Value = "21"

$ws.Range("D8").Value = "1.396"
$ws.Range("E8").Value = "2.05%"
$ws.Range("G8").Value = "21"

$ws.Range("D9").Value = "0.9069"
$ws.Range("E9").Value = "-3.66%"
$ws.Range("G9").Value = "21"

$ws.Range("D10").Value = "0.1612"
$ws.Range("E10").Value = "2.78%"
$ws.Range("G10").Value = "21"

$ws.Range("D11").Value = "0.07678"
$ws.Range("E11").Value = "12.70%"
$ws.Range("G11").Value = "21"

$ws.Range("D12").Value = "0.07721"
$ws.Range("E12").Value = "1.67%"
$ws.Range("G12").Value = "21"

$ws.Range("D13").Value = "0.02913"
$ws.Range("E13").Value = "-1.29%"
$ws.Range("G13").Value = "21"

$ws.Range("D14").Value = "0.09031"
$ws.Range("E14").Value = "0.25%"
$ws.Range("G14").Value = "21"

$ws.Range("D15").Value = "0.001599"
$ws.Range("E15").Value = "0.57%"
$ws.Range("G15").Value = "21"

$ws.Range("D16").Value = "0.0006512"
$ws.Range("E16").Value = "0.65%"
$ws.Range("G16").Value = "21"

$ws.Range("D17").Value = "0.006092"
$ws.Range("E17").Value = "-3.57%"
$ws.Range("G17").Value = "21"

$ws.Range("D18").Value = "3.493"
$ws.Range("E18").Value = "1.26%"
$ws.Range("G18").Value = "21"

$ws.Range("E19").Value = "-0.90%"
$ws.Range("G19").Value = "21"

$ws.Range("D20").Value = "0.3231"
$ws.Range("E20").Value = "0.67%"
$ws.Range("G20").Value = "21"

$ws.Range("D21").Value = "0.1348"
$ws.Range("E21").Value = "2.88%"
$ws.Range("G21").Value = "21"

$ws.Range("D22").Value = "4.012"
$ws.Range("E22").Value = "-1.21%"
$ws.Range("G22").Value = "21"

$ws.Range("D23").Value = "0.1600"
$ws.Range("E23").Value = "3.26%"
$ws.Range("G23").Value = "21"

$ws.Range("D24").Value = "0.04530"
$ws.Range("E24").Value = "1.09%"
$ws.Range("G24").Value = "21"

$ws.Range("D25").Value = "0.001213"
$ws.Range("E25").Value = "2.94%"
$ws.Range("G25").Value = "21"

$ws.Range("D26").Value = "0.004161"
$ws.Range("E26").Value = "-7.41%"
$ws.Range("G26").Value = "21"

$ws.Range("D27").Value = "0.0001170"
$ws.Range("E27").Value = "-6.15%"
$ws.Range("G27").Value = "21"

$ws.Range("D28").Value = "0.0001670"
$ws.Range("E28").Value = "3.41%"
$ws.Range("G28").Value = "21"

$ws.Range("G29").Value = "21"

$ws.Range("G30").Value = "21"

$ws.Range("G31").Value = "21"

$ws.Range("G32").Value = "21"

$ws.Range("G33").Value = "21"

$ws.Range("G34").Value = "21"

$ws.Range("G35").Value = "21"

$ws.Range("G36").Value = "21"

$ws.Range("G37").Value = "21"

$ws.Range("G38").Value = "21"

$ws.Range("G39").Value = "21"

$ws.Range("D40").Value = "0.04385"
$ws.Range("E40").Value = "4.36%"
$ws.Range("G40").Value = "21"

$ws.Range("D41").Value = "0.007003"
$ws.Range("E41").Value = "4.13%"
$ws.Range("G41").Value = "21"

$ws.Range("D42").Value = "0.1254"
$ws.Range("E42").Value = "-0.20%"
$ws.Range("G42").Value = "21"

$ws.Range("D43").Value = "0.002071"
$ws.Range("E43").Value = "2.75%"
$ws.Range("G43").Value = "21"

$ws.Range("D44").Value = "0.01185"
$ws.Range("E44").Value = "-3.59%"
$ws.Range("G44").Value = "21"

$ws.Range("D45").Value = "0.00005851"
$ws.Range("E45").Value = "3.63%"
$ws.Range("G45").Value = "21"

$ws.Range("G46").Value = "21"

$ws.Range("D47").Value = "0.01300"
$ws.Range("E47").Value = "-0.29%"
$ws.Range("G47").Value = "21"

$ws.Range("G48").Value = "21"

$ws.Range("G49").Value = "21"

$ws.Range("G50").Value = "21"

$ws.Range("G51").Value = "21"

